$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert a new row before row 310 for the new "TableCellCollection"
#    class entry (Word > TableCellCollection > class).
# ------------------------------------------------------------------
$ws.Rows.Item(310).Insert()

$ws.Range("A310").Value = "Word"
$ws.Range("B310").Value = "TableCellCollection"
$ws.Range("D310").Value = "class"
$ws.Range("E310").Value = "word-tables-manage-formatting"
$ws.Range("F310").Value = "getTableCellAlignment"

# Match formatting to the neighboring rows in this block (font-only
# style on A/B, vertical-centered style on E/F), and give the empty
# C310 cell the same "blank but formatted" treatment as the rest of
# the block.
$ws.Range("A311:F311").Copy()
$ws.Range("A310:F310").PasteSpecial(-4122)

$ws.Range("A310").Value = "Word"
$ws.Range("B310").Value = "TableCellCollection"
$ws.Range("C310").Value = ""
$ws.Range("D310").Value = "class"
$ws.Range("E310").Value = "word-tables-manage-formatting"
$ws.Range("F310").Value = "getTableCellAlignment"

# ------------------------------------------------------------------
# 2) Insert a new row before (the now shifted) row 315 for the new
#    "cells" member on TableRow (Word > TableRow > cells).
# ------------------------------------------------------------------
$ws.Rows.Item(315).Insert()

$ws.Range("A314:F314").Copy()
$ws.Range("A315:F315").PasteSpecial(-4122)

$ws.Range("A315").Value = "Word"
$ws.Range("B315").Value = "TableRow"
$ws.Range("C315").Value = "cells"
$ws.Range("D315").Value = ""
$ws.Range("E315").Value = "word-tables-manage-formatting"
$ws.Range("F315").Value = "getTableCellAlignment"
$ws.Range("D315").ClearContents()

# ------------------------------------------------------------------
# 3) Grow the "Snippets" table to cover the two newly inserted rows.
# ------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lastRow = $ws.UsedRange.Rows.Count
$lo.Resize($ws.Range("A1:F" + $lastRow))

# ------------------------------------------------------------------
# 4) Restore the view: scroll back to the top of the frozen pane and
#    select A2 (matches the pre-edit-session view state).
# ------------------------------------------------------------------
$ws.Range("A2").Select()
